$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Set "Runmode" to "N" for suites C, D, E, F (rows 4-7) so they no longer
# run, while suites A and B (rows 2-3) keep running ("Y"), i.e. "Running
# suites A and B parallely".
$ws.Range("C4").Value = "N"
$ws.Range("C5").Value = "N"
$ws.Range("C6").Value = "N"
$ws.Range("C7").Value = "N"
